$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (knn)
$ws.Range("B4").Value = 0.6130949948175259
$ws.Range("C4").Value = 0.609
$ws.Range("D4").Value = 0.6308209611002592
$ws.Range("E4").Value = 0.6210000000000001
$ws.Range("F4").Value = 0.5048500116638125
$ws.Range("G4").Value = 0.5000000000000001
$ws.Range("H4").Value = 0.5131360300703008
$ws.Range("I4").Value = 0.5130000000000001
$ws.Range("J4").Value = 0.618624316015203
$ws.Range("K4").Value = 0.6469999999999999
$ws.Range("L4").Value = 0.5980672807509251
$ws.Range("M4").Value = 0.601

# Row 5 (nb)
$ws.Range("B5").Value = 0.3049745498958462
$ws.Range("C5").Value = 0.229
$ws.Range("D5").Value = 0.5935315452091767
$ws.Range("E5").Value = 0.5645
$ws.Range("J5").Value = 0.6013616104472984
$ws.Range("K5").Value = 0.7590000000000001
$ws.Range("L5").Value = 0.5635485265153199
$ws.Range("M5").Value = 0.5685

# Row 6 (rf)
$ws.Range("B6").Value = 0.6800672834381298
$ws.Range("C6").Value = 0.696
$ws.Range("D6").Value = 0.6789847778345015
$ws.Range("E6").Value = 0.6855
$ws.Range("F6").Value = 0.5105136622013136
$ws.Range("G6").Value = 0.512
$ws.Range("H6").Value = 0.5109260796769836
$ws.Range("I6").Value = 0.5125000000000001
$ws.Range("J6").Value = 0.7006441631166135
$ws.Range("K6").Value = 0.6830000000000001
$ws.Range("L6").Value = 0.7362084704971299
$ws.Range("M6").Value = 0.7155
